$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates as Excel serial numbers, matching column A's existing date format)
$newRows = @(
    @{ Row = 234; Date = 44308; B = 2; C = 18; D = 159.4472495349455 },
    @{ Row = 235; Date = 44309; B = 2; C = 14; D = 124.0145274160687 },
    @{ Row = 236; Date = 44310; B = 7; C = 17; D = 150.5890690052263 },
    @{ Row = 237; Date = 44311; B = 1; C = 16; D = 141.7308884755071 },
    @{ Row = 238; Date = 44312; B = 2; C = 18; D = 159.4472495349455 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Copy formatting (style/border/alignment/number format) from the row above, column A
    $srcA = $ws.Cells.Item($rowNum - 1, 1)
    $dstA = $ws.Cells.Item($rowNum, 1)
    $srcA.Copy($dstA)
    $dstA.Value = $r.Date

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
}
